# edit.ps1 - apply the changes described by the target diff:
#   1) Update the cached "datetimeFigureOut" date field text from
#      16/12/2022 -> 17/12/2022 on the slide master and every slide
#      layout's "Date Placeholder" shape.
#   2) On slide 9, split the "25875 total rows of data" run into two
#      runs so the text reads "8625 total rows of data" (new leading
#      run "8625 " + the original trailing run "total rows of data").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Date placeholders (slide master + every custom layout)
# ---------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "16/12/2022") {
                $shp.TextFrame.TextRange.Text = "17/12/2022"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------
# 2) Slide 9 - "25875 total rows of data" -> "8625 total rows of data"
#    (split into two runs: "8625 " + "total rows of data")
# ---------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
for ($i = 1; $i -le $slide9.Shapes.Count; $i++) {
    $shp = $slide9.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "25875 total rows of data") {
        $tr = $shp.TextFrame.TextRange
        $prefix = $tr.Characters(1, 6)
        $prefix.Text = "8625 "
    }
}
